# Fruta / hortaliza, semanal
# Insert two new weekly records at rows 38-39 (pushing the existing
# rows 38-51 down to 40-53) in the "Durazno" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 38, shifting the
# existing data (old rows 38-51) down to rows 40-53.
$ws.Range("A38:A39").EntireRow.Insert()

# --- New row 38: Phillips Cling / Segunda ---
$ws.Range("A38").Value = 1
$ws.Range("B38").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C38").Value = "Arica y Parinacota"
$ws.Range("D38").Value = 44637
$ws.Range("E38").Value = 15
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100103
$ws.Range("H38").Value = "Frutos de hueso (carozo)"
$ws.Range("I38").Value = 100103004
$ws.Range("J38").Value = "Durazno"
$ws.Range("K38").Value = "Phillips Cling"
$ws.Range("L38").Value = "Segunda"
$ws.Range("M38").Value = 300
$ws.Range("N38").Value = 18000
$ws.Range("O38").Value = 20000
$ws.Range("P38").Value = 19000
$ws.Range("Q38").Value = "$/bandeja 18 kilos granel"
$ws.Range("R38").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S38").Value = 1056
$ws.Range("T38").Value = 18

# --- New row 39: September Snow / Segunda ---
$ws.Range("A39").Value = 1
$ws.Range("B39").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C39").Value = "Arica y Parinacota"
$ws.Range("D39").Value = 44637
$ws.Range("E39").Value = 15
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100103
$ws.Range("H39").Value = "Frutos de hueso (carozo)"
$ws.Range("I39").Value = 100103004
$ws.Range("J39").Value = "Durazno"
$ws.Range("K39").Value = "September Snow"
$ws.Range("L39").Value = "Segunda"
$ws.Range("M39").Value = 250
$ws.Range("N39").Value = 19000
$ws.Range("O39").Value = 20000
$ws.Range("P39").Value = 19500
$ws.Range("Q39").Value = "$/bandeja 18 kilos granel"
$ws.Range("R39").Value = "Región de O'Higgins"
$ws.Range("S39").Value = 1083
$ws.Range("T39").Value = 18
